# Weekly update: insert a new price-report row for "Femacal de La Calera -
# Perejil" at row 9 (pushing the existing rows 9-31 down to rows 10-32),
# then populate the newly inserted row with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9 - Excel shifts rows 9:31 down to 10:32,
# carrying their formatting (incl. the date-style D column) with them.
$ws.Rows.Item(9).Insert()

# Fill the newly inserted row 9 with the new week's data.
$ws.Cells.Item(9, 1).Value  = 3
$ws.Cells.Item(9, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(9, 3).Value  = "Coquimbo"
$ws.Cells.Item(9, 4).Value  = 44965
$ws.Cells.Item(9, 5).Value  = 5
$ws.Cells.Item(9, 6).Value  = 100112044
$ws.Cells.Item(9, 7).Value  = "Perejil"
$ws.Cells.Item(9, 8).Value  = "Sin especificar"
$ws.Cells.Item(9, 9).Value  = "Primera"
$ws.Cells.Item(9, 10).Value = 87
$ws.Cells.Item(9, 11).Value = 3000
$ws.Cells.Item(9, 12).Value = 3000
$ws.Cells.Item(9, 13).Value = 3000
$ws.Cells.Item(9, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(9, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(9, 16).Value = 1000
$ws.Cells.Item(9, 17).Value = 3
$ws.Cells.Item(9, 18).Value = "Hortaliza"
